$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for "Femacal de La Calera" / Acelga is inserted
# at row 279, pushing the existing rows 279:389 down to 280:390.
$ws.Rows(279).Insert()

$ws.Range("A279").Value = 3
$ws.Range("B279").Value = "Femacal de La Calera"
$ws.Range("C279").Value = "Coquimbo"
$ws.Range("D279").Value = 44784
$ws.Range("E279").Value = 5
$ws.Range("F279").Value = 100112009
$ws.Range("G279").Value = "Acelga"
$ws.Range("H279").Value = "Sin especificar"
$ws.Range("I279").Value = "Primera"
$ws.Range("J279").Value = 220
$ws.Range("K279").Value = 3300
$ws.Range("L279").Value = 3500
$ws.Range("M279").Value = 3400
$ws.Range("N279").Value = "$/docena de atados (6 kilos)"
$ws.Range("O279").Value = "Provincia de Quillota"
$ws.Range("P279").Value = 567
$ws.Range("Q279").Value = 6
$ws.Range("R279").Value = "Hortaliza"
